$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column O (15th column) of data for year 2021, matching column N formatting
$ws.Range("N3:N5").Copy()
$ws.Range("O3:O5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("O3").Value = 2021
$ws.Range("O4").Value = 14
$ws.Range("O5").Value = 1252.8

$ws.Range("O9").Select()
